# Sprint 1 burndown update:
#  - Fill in "Actual" hours (column I) for the Sprint 1 retrospective rows on
#    the Data sheet (previously blank, now logged).
#  - Update the Actual burndown value for the last data point (P7).
#  - Move the selection on the Data sheet from P7 to P8.
#
# NOTE: selecting a range on a worksheet that is not already the active sheet
# implicitly activates it (and thus would flip which sheet tab is marked
# "selected" in the saved file). The source workbook has "Burndown" as the
# active/selected tab, and that does not change in the target, so we restore
# it as the active sheet after moving the selection on "Data".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("I3").Value = 0
$ws.Range("I4").Value = 2
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 0
$ws.Range("I7").Value = 0

$ws.Range("P7").Value = 8

$ws.Range("P8").Select() | Out-Null

$wb.Worksheets.Item("Burndown").Activate() | Out-Null
